# Convert NameConversions to use the new list of syntax
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Underscores")

# Update the "List_Property(0) of" label to the new "List_Property list of" syntax
$ws.Range("C6").Value = "List_Property list of"

# Update the nested "With Properties" label under the list to "With Item"
$ws.Range("D7").Value = "With Item"

# Move the selection / active cell on this sheet
$ws.Activate()
$ws.Range("D8").Select()
